# Budget for Sunlight Sensor - add a "Shipping Cost (Canada Post)" line to
# the budget table (row 19) and roll it into the Grand Total formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: new "Shipping Cost (Canada Post)" line -----------------------
# A19 gets the same plain "label" formatting used by the other section
# separators in column A (e.g. A7, A10, A17).
$ws.Range("A19").Font.Bold = $false

# F19: bold + centered (matches the "Subtotal"-style emphasis used for the
# other bold/centered cells in the sheet).
$ws.Range("F19").Font.Bold = $true
$ws.Range("F19").HorizontalAlignment = -4108   # xlCenter

# G19: centered, regular weight.
$ws.Range("G19").HorizontalAlignment = -4108   # xlCenter

# H19: the new shipping-cost label, centered + wrapped (mirrors the H5
# header "Shipping Fee").
$ws.Range("H19").Value = "Shipping Cost (Canada Post)"
$ws.Range("H19").HorizontalAlignment = -4108   # xlCenter
$ws.Range("H19").WrapText = $true

# I19: the shipping cost itself, formatted like the other CAD$ cost cells
# (I6 / I10 / I14 / I20).
$ws.Range("I19").NumberFormat = $ws.Range("I6").NumberFormat
$ws.Range("I19").Value = 15.41

# Row 19 grows to accommodate the wrapped label text.
$ws.Rows.Item(19).RowHeight = 60

# --- Grand total formula now also includes the new shipping line ----------
$ws.Range("I20").Formula = "=SUM(I14:I15,I10,I6,I19)"

# --- Columns H and I share the same width now -----------------------------
$ws.Columns.Item(9).ColumnWidth = $ws.Columns.Item(8).ColumnWidth

# --- View state: selection moved to H24, scrolled down a bit --------------
$ws.Range("H24").Select()

# --- Print scaling tweak ----------------------------------------------------
$ws.PageSetup.Zoom = 79
